# Auto-generated edit script applying the Hyperion_Profits market-data refresh
# across sheets ALC, ARM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 138.5
$ws.Range("I9").Value = 158.33333
$ws.Range("K9").Value = 158.33333
$ws.Range("M9").Value = 10.66667000000001

# Row 43
$ws.Range("H43").Value = 26317088
$ws.Range("J43").Value = 1549.2858
$ws.Range("L43").Value = 1549.2858
$ws.Range("N43").Value = -1687.2858

# Row 129
$ws.Range("H129").Value = 1585
$ws.Range("I129").Value = 1585
$ws.Range("K129").Value = 4755
$ws.Range("M129").Value = 245

# Row 130
$ws.Range("H130").Value = 398998
$ws.Range("J130").Value = 398998
$ws.Range("L130").Value = 398998
$ws.Range("N130").Value = -409038

# Row 132
$ws.Range("H132").Value = 43481372
$ws.Range("I132").Value = 55558256
$ws.Range("K132").Value = 166674768
$ws.Range("M132").Value = -166672238

# Row 135
$ws.Range("H135").Value = 1121.5
$ws.Range("I135").Value = 627.8823
$ws.Range("K135").Value = 5650.9407
$ws.Range("M135").Value = -3115.9407

# Row 137
$ws.Range("H137").Value = 179357.8
$ws.Range("J137").Value = 2290.25
$ws.Range("L137").Value = 6870.75
$ws.Range("N137").Value = -11970.75

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 9999
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 9999
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 9999
$ws.Range("M6").ClearContents() # was 73
$ws.Range("N6").Value = -10345

# Row 9
$ws.Range("H9").Value = 37999.668
$ws.Range("J9").Value = 37999.668
$ws.Range("L9").Value = 37999.668
$ws.Range("N9").Value = -38339.668

# Row 20
$ws.Range("H20").Value = 37999.668
$ws.Range("J20").Value = 37999.668
$ws.Range("L20").Value = 37999.668
$ws.Range("N20").Value = -38539.668

# Row 32
$ws.Range("H32").Value = 4165.1016
$ws.Range("I32").Value = 2798.652
$ws.Range("J32").Value = 9000.23
$ws.Range("K32").Value = 2798.652
$ws.Range("L32").Value = 9000.23
$ws.Range("M32").Value = -2511.652
$ws.Range("N32").Value = -9574.23

# Row 61
$ws.Range("H61").Value = 8164.8335
$ws.Range("I61").Value = 10248.5
$ws.Range("K61").Value = 10248.5
$ws.Range("M61").Value = -10036.5

# Row 74
$ws.Range("H74").Value = 72962.88
$ws.Range("J74").Value = 272441.7
$ws.Range("L74").Value = 272441.7
$ws.Range("N74").Value = -274189.7

# Row 77
$ws.Range("H77").Value = 72962.88
$ws.Range("J77").Value = 272441.7
$ws.Range("L77").Value = 1362208.5
$ws.Range("N77").Value = -1370944.5

# Row 132
$ws.Range("H132").Value = 3500.3333
$ws.Range("J132").Value = 6230.5
$ws.Range("L132").Value = 18691.5
$ws.Range("N132").Value = -23751.5

# Row 136
$ws.Range("H136").Value = 8164.8335
$ws.Range("I136").Value = 10248.5
$ws.Range("K136").Value = 30745.5
$ws.Range("M136").Value = -28195.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 15537.177
$ws.Range("I31").Value = 2951
$ws.Range("J31").Value = 16666.705
$ws.Range("K31").Value = 2951
$ws.Range("L31").Value = 16666.705
$ws.Range("M31").Value = -2656
$ws.Range("N31").Value = -17256.705

# Row 34
$ws.Range("H34").Value = 15537.177
$ws.Range("I34").Value = 2951
$ws.Range("J34").Value = 16666.705
$ws.Range("K34").Value = 2951
$ws.Range("L34").Value = 16666.705
$ws.Range("M34").Value = -2749
$ws.Range("N34").Value = -17070.705

# Row 132
$ws.Range("H132").Value = 79915.46000000001
$ws.Range("I132").Value = 93536.45
$ws.Range("K132").Value = 280609.35
$ws.Range("M132").Value = -278079.35

# Row 134
$ws.Range("H134").Value = 2395.4187
$ws.Range("I134").Value = 1642.3793
$ws.Range("K134").Value = 4927.1379
$ws.Range("M134").Value = -2392.1379

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 15259.571
$ws.Range("I5").Value = 1350
$ws.Range("K5").Value = 4050
$ws.Range("M5").Value = -3938

# Row 16
$ws.Range("H16").Value = 467.25
$ws.Range("J16").Value = 467.25
$ws.Range("L16").Value = 1401.75
$ws.Range("N16").Value = -1747.75

# Row 68
$ws.Range("H68").Value = 2408.3125
$ws.Range("I68").Value = 2423.4443
$ws.Range("J68").Value = 2402.3914
$ws.Range("K68").Value = 7270.3329
$ws.Range("L68").Value = 7207.174199999999
$ws.Range("M68").Value = -6459.3329
$ws.Range("N68").Value = -8829.174199999999

# Row 71
$ws.Range("H71").Value = 2408.3125
$ws.Range("I71").Value = 2423.4443
$ws.Range("J71").Value = 2402.3914
$ws.Range("K71").Value = 21810.9987
$ws.Range("L71").Value = 21621.5226
$ws.Range("M71").Value = -17754.9987
$ws.Range("N71").Value = -29733.5226

# Row 135
$ws.Range("H135").Value = 15259.571
$ws.Range("I135").Value = 1350
$ws.Range("K135").Value = 12150
$ws.Range("M135").Value = -9615

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 4996.5
$ws.Range("I11").Value = 4996.5
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 4996.5
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -4857.5
$ws.Range("N11").ClearContents() # was -15278

# Row 14
$ws.Range("H14").Value = 977.5
$ws.Range("I14").Value = 977.5
$ws.Range("K14").Value = 977.5
$ws.Range("M14").Value = -809.5

# Row 80
$ws.Range("H80").Value = 65550610
$ws.Range("J80").Value = 2452.75
$ws.Range("L80").Value = 2452.75
$ws.Range("N80").Value = -4448.75

# Row 83
$ws.Range("H83").Value = 65550610
$ws.Range("J83").Value = 2452.75
$ws.Range("L83").Value = 12263.75
$ws.Range("N83").Value = -22247.75

# Row 126
$ws.Range("H126").Value = 4774096.5
$ws.Range("I126").Value = 2529612
$ws.Range("J126").Value = 9263066
$ws.Range("K126").Value = 7588836
$ws.Range("L126").Value = 27789198
$ws.Range("M126").Value = -7586366
$ws.Range("N126").Value = -27794138

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3984.0527
$ws.Range("I7").Value = 2594.2666
$ws.Range("K7").Value = 2594.2666
$ws.Range("M7").Value = -2482.2666

# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents() # was -2828

# Row 22
$ws.Range("H22").Value = 64239.43
$ws.Range("J22").Value = 914.75
$ws.Range("L22").Value = 914.75
$ws.Range("N22").Value = -1504.75

# Row 27
$ws.Range("H27").Value = 64239.43
$ws.Range("J27").Value = 914.75
$ws.Range("L27").Value = 914.75
$ws.Range("N27").Value = -1128.75

# Row 41
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents() # was -15876

# Row 46
$ws.Range("H46").Value = 3108699.8
$ws.Range("I46").Value = 7246999.5
$ws.Range("J46").Value = 4974.75
$ws.Range("K46").Value = 7246999.5
$ws.Range("L46").Value = 4974.75
$ws.Range("M46").Value = -7246811.5
$ws.Range("N46").Value = -5350.75

# Row 126
$ws.Range("H126").Value = 3984.0527
$ws.Range("I126").Value = 2594.2666
$ws.Range("K126").Value = 7782.7998
$ws.Range("M126").Value = -5312.7998

# Row 136
$ws.Range("H136").Value = 37454.367
$ws.Range("I136").Value = 44818.168
$ws.Range("J136").Value = 7999.1665
$ws.Range("K136").Value = 134454.504
$ws.Range("L136").Value = 23997.4995
$ws.Range("M136").Value = -131904.504
$ws.Range("N136").Value = -29097.4995

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 18936.375
$ws.Range("I62").Value = 37163.668
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 37163.668
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -36539.668
$ws.Range("N62").Value = -9248

# Row 65
$ws.Range("H65").Value = 18936.375
$ws.Range("I65").Value = 37163.668
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 185818.34
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -182698.34
$ws.Range("N65").Value = -46240

# Row 122
$ws.Range("H122").Value = 3473.5356
$ws.Range("I122").Value = 3209.7727
$ws.Range("K122").Value = 9629.3181
$ws.Range("M122").Value = -7179.3181

# Row 136
$ws.Range("H136").Value = 5060.3945
$ws.Range("I136").Value = 5928.3477
$ws.Range("J136").Value = 3729.5334
$ws.Range("K136").Value = 17785.0431
$ws.Range("L136").Value = 11188.6002
$ws.Range("M136").Value = -15235.0431
